$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = 15000000.0
$ws.Range("C7").Value = 14000000.0
$ws.Range("D7").Value = 17891000.0
$ws.Range("E7").Value = 18852000.0
$ws.Range("F7").Value = 19176000.0
